$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "e"
$ws.Range("C2").Value = "e"
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 48.26

# Row 3
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "e"
$ws.Range("C3").Value = "e"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 16.09

# Row 4
$ws.Range("A4").Value = 4
